$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.196.82"
$ws.Range("E2").Value = "  +1.55%  "

# Row 3
$ws.Range("D3").Value = "1.782.51"
$ws.Range("E3").Value = "  +0.55%  "

# Row 4
$ws.Range("E4").Value = "  +0.23%  "

# Row 5
$ws.Range("D5").Value = "'225.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "

# Row 6
$ws.Range("E6").Value = "  +0.49%  "

# Row 7
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("D8").Value = "'31.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.33%  "

# Row 9
$ws.Range("E9").Value = "  +0.92%  "

# Row 10
$ws.Range("E10").Value = "  +0.35%  "

# Row 11
$ws.Range("D11").Value = "'0.0946"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.22%  "

# Row 12
$ws.Range("D12").Value = "2.039.96"
$ws.Range("E12").Value = "  +0.73%  "

# Row 13
$ws.Range("D13").Value = "'10.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.78%  "

# Row 14
$ws.Range("D14").Value = "1.771.42"
$ws.Range("E14").Value = "  -0.12%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.623"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.43%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.095.26"
$ws.Range("E16").Value = "  +1.25%  "

# Row 17
$ws.Range("E17").Value = "  +1.47%  "

# Row 18
$ws.Range("D18").Value = "'68.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.30%  "

# Row 19
$ws.Range("D19").Value = "'245.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.33%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0787"
$ws.Range("E20").Value = "  +1.53%  "

# Row 21
$ws.Range("D21").Value = "'10.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.82%  "

# Row 22
$ws.Range("E22").Value = "  +0.17%  "

# Row 23
$ws.Range("E23").Value = "  +2.57%  "

# Row 24
$ws.Range("E24").Value = "  -0.27%  "

# Row 25
$ws.Range("D25").Value = "'161.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "

# Row 26
$ws.Range("E26").Value = "  +2.52%  "

# Row 27
$ws.Range("E27").Value = "  +1.45%  "

# Row 28
$ws.Range("E28").Value = "  +1.74%  "

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.29%  "

# Row 30
$ws.Range("E30").Value = "  +0.77%  "

# Row 31
$ws.Range("D31").Value = "'0.0519"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.65%  "

# Row 32
$ws.Range("E32").Value = "  +3.97%  "

# Row 33
$ws.Range("D33").Value = "'3.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.10%  "

# Row 34
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("D35").Value = "1.442.81"
$ws.Range("E35").Value = "  +4.59%  "

# Row 36
$ws.Range("D36").Value = "'0.654"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.11%  "

# Row 37
$ws.Range("D37").Value = "'2.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.71%  "

# Row 38
$ws.Range("E38").Value = "  +3.88%  "

# Row 39
$ws.Range("E39").Value = "  +1.51%  "

# Row 40
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").Value = "'2.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.75%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'80.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.75%  "

# Row 42
$ws.Range("D42").Value = "'0.924"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.01%  "

# Row 43
$ws.Range("E43").Value = "  +0.72%  "

# Row 44
$ws.Range("D44").Value = "'13.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.69%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0139"
$ws.Range("E45").Value = "  +2.77%  "

# Row 46
$ws.Range("D46").Value = "'0.0509"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.92%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'6.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.97%  "

# Row 48
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.17%  "

# Row 49
$ws.Range("D49").Value = "1.942.09"
$ws.Range("E49").Value = "  +1.02%  "

# Row 50
$ws.Range("D50").Value = "'105.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "

# Row 51
$ws.Range("E51").Value = "  +0.18%  "
